$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-7 (columns D, L, M, N, O, P, Q, R, S, T) are cyclically
# rotated "up" by one row: the original row 3 values move into row 2, the
# original row 4 values move into row 3, ..., the original row 7 values move
# into row 6, and the original row 2 values wrap around into row 7.
# Columns A, B, C, E, F, G, H, I, J, K are identical across all rows and stay
# unchanged.

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Capture the original values for rows 2..7 before making any changes.
$original = @{}
for ($r = 2; $r -le 7; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Row r (2..6) gets the original values of row r+1; row 7 gets the original
# values of row 2.
for ($r = 2; $r -le 7; $r++) {
    $srcRow = $r + 1
    if ($srcRow -gt 7) { $srcRow = 2 }
    $src = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $src[$c]
    }
}
